# Updates cryptos list values (prices + % volume) per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.069.07"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.309.66"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.62"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.99"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +6.00%  "

$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +5.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.68"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +7.31%  "

$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.117"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +3.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.98"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +15.64%  "

$ws.Range("E14").Value = "  +3.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.687.70"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.306.49"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("E17").Value = "  +3.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.005.63"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +2.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +8.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +3.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +1.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.98"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.66"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("E24").Value = "  +12.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +5.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.72"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.61"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +2.64%  "

$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.04"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.08"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +3.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0696"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.46%  "

$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.50%  "

$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0289"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +4.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.29"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +7.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.59"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("E47").Value = "  +3.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.08"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +7.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.529.33"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.53"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +3.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.58"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.03%  "

# Row 42/43: ApeXProtocol and Maker swapped positions, with new price/volume data
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.001.49"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +1.39%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.30"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -5.04%  "
